$d = $word.ActiveDocument
$r = $d.Content
$xml = $r.WordOpenXML

# The shape's spPr contains an <a:ln> (outline/line) element that isn't
# needed for this test; remove it entirely (it was between </a:gradFill>
# and </wps:spPr>).
$pattern = '(?s)<a:ln cap="rnd" w="57240">.*?</a:ln>'
$matches = [System.Text.RegularExpressions.Regex]::Matches($xml, $pattern)
if ($matches.Count -ne 1) {
    throw "expected exactly one <a:ln> block to remove, found $($matches.Count)"
}

$newXml = [System.Text.RegularExpressions.Regex]::Replace($xml, $pattern, '')
$r.InsertXML($newXml)
